$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Clear()

$ws.Range("A1").Value = "Reference"
$ws.Range("B1").Value = "QID"
$ws.Range("C1").Value = "Morphology"
$ws.Range("D1").Value = "Alias_1"
$ws.Range("E1").Value = "Alias_2"
$ws.Range("F1").Value = "Alias_3"
$ws.Range("G1").Value = "Locked (Y/N)"
$ws.Range("H1").Value = "Primary Alias"

$ws.Range("A2").Value = "10.3389/fchem.2018.00407"
$ws.Range("B2").Value = "Q11086567"
$ws.Range("C2").Value = "Fragment"
$ws.Range("D2").Value = "Fragments"
$ws.Range("G2").Value = "Y"

$ws.Range("A3").Value = "10.3389/fchem.2018.00407"
$ws.Range("B3").Value = "Q161"
$ws.Range("C3").Value = "Fiber"
$ws.Range("D3").Value = "Fibers"
$ws.Range("E3").Value = "Fibres"
$ws.Range("G3").Value = "N"

$ws.Range("A4").Value = "10.1016/j.watres.2019.02.054"
$ws.Range("B4").Value = "Q109875324"
$ws.Range("C4").Value = "Nurdle"
$ws.Range("D4").Value = "Pellet"
$ws.Range("E4").Value = "Pellets"
$ws.Range("G4").Value = "N"
$ws.Range("H4").Value = "Nurdle"

$ws.Range("A5").Value = "10.3389/fchem.2018.00407"
$ws.Range("B5").Value = "Q1137203"
$ws.Range("C5").Value = "Film"
$ws.Range("D5").Value = "Films"
$ws.Range("G5").Value = "N"
$ws.Range("H5").Value = "Thin Film"

$ws.Range("A6").Value = "10.3389/fchem.2018.00407"
$ws.Range("B6").Value = "Q215414"
$ws.Range("C6").Value = "Foam"
$ws.Range("G6").Value = "N"

$ws.Range("A7").Value = "10.1016/j.watres.2019.02.054"
$ws.Range("B7").Value = "Q12507"
$ws.Range("C7").Value = "Sphere"
$ws.Range("D7").Value = "Spheres"
$ws.Range("E7").Value = "Spherical"
$ws.Range("F7").Value = "Sphericals"
$ws.Range("G7").Value = "Y"

$ws.Range("A8").Value = "10.1016/j.watres.2019.02.054"
$ws.Range("B8").Value = "Q37105"
$ws.Range("C8").Value = "Line"
$ws.Range("G8").Value = "N"

$ws.Range("A9").Value = "10.1016/j.watres.2019.02.054"
$ws.Range("B9").Value = "Q1053956"
$ws.Range("C9").Value = "Bead"
$ws.Range("G9").Value = "N"

$ws.Range("A10").Value = "10.1016/j.watres.2019.02.054"
$ws.Range("B10").Value = "Q66539740"
$ws.Range("C10").Value = "Sheets"
$ws.Range("D10").Value = "Flake"
$ws.Range("G10").Value = "N"
$ws.Range("H10").Value = "Sheet"

$ws.Range("A11").Value = "10.1007/s11783-021-1492-5"
$ws.Range("C11").Value = "Films/fragments"

$ws.Range("A11").Font.Color = 0

$ws.Range("A11").Select()